# Applies the German copy-edit pass described in the commit diff:
#  - Slide 2: "Agenda" -> "Tagesordnung" (title) and the agenda bullet
#             "Produktbeschreibung" -> "Produktbeschreibung:"
#  - Slide 4: Title "Produktbeschreibung" -> "Produktbeschreibung:"
#  - Slide 5: two table-cell bullet rewrites (Authentischer Blend.../
#             Gesundheitsfoerdernde Inhaltsstoffe...)
#  - Slide 6: five table-cell bullet rewrites (Vielseitige.../Nachhaltig
#             gewonnen.../Elegante Verpackung.../Mit Kundenzufriedenheits.../
#             Perfekt geeignet...)
#
# Each bullet in the feature tables is stored as several <a:r> runs that
# share identical formatting (the visible sentence, a literal space, and a
# trailing sentence). PowerPoint's object model folds consecutive
# same-formatted runs into a single logical Run, so indexing Runs(1,1) gets
# only the *first* underlying XML run - exactly the span the diff rewrites -
# leaving the later runs (and their text) untouched.

$p = $ppt.ActivePresentation

# ---- Slide 2 -------------------------------------------------------------
$s2 = $p.Slides.Item(2)

# Title: "Agenda" -> "Tagesordnung"
$s2.Shapes.Item("Title 1").TextFrame.TextRange.Runs(1,1).Text = "Tagesordnung"

# Agenda list bullet (2nd paragraph): "Produktbeschreibung" -> "Produktbeschreibung:"
$contentPh = $s2.Shapes.Item("Content Placeholder 2").TextFrame.TextRange
$contentPh.Paragraphs(2,1).Runs(1,1).Text = "Produktbeschreibung:"

# ---- Slide 4 ---------------------------------------------------------------
$s4 = $p.Slides.Item(4)

# Title: "Produktbeschreibung" -> "Produktbeschreibung:"
$s4.Shapes.Item("Title 1").TextFrame.TextRange.Runs(1,1).Text = "Produktbeschreibung:"

# ---- Slide 5 ---------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$tbl5 = $s5.Shapes.Item("Content Placeholder 4").Table

$tbl5.Cell(4,1).Shape.TextFrame.TextRange.Paragraphs(1,1).Runs(1,1).Text = `
    "Authentische Mischung: Unser Chai ist eine harmonische Mischung aus hochwertigen Schwarzteeblättern und einer charakteristischen Auswahl an gemahlenen Gewürzen wie Zimt, Kardamom, Nelken, Ingwer und schwarzem Pfeffer."

$tbl5.Cell(4,2).Shape.TextFrame.TextRange.Paragraphs(1,1).Runs(1,1).Text = `
    "Gesundheitsfördernde Zutaten: Alle Inhaltsstoffe des Mystic Spice Chai Tea werden aufgrund ihrer natürlichen gesundheitsfördernden Eigenschaften ausgewählt."

# ---- Slide 6 ---------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$tbl6 = $s6.Shapes.Item("Content Placeholder 4").Table

$tbl6.Cell(2,2).Shape.TextFrame.TextRange.Paragraphs(1,1).Runs(1,1).Text = `
    "Vielfältige Zubereitungsmöglichkeiten: Ob Sie Ihren Chai dampfend heiß, als erfrischenden Eistee oder als cremigen Latte mögen, unsere Mischung ist vielseitig, um allen Vorlieben gerecht zu werden."

$tbl6.Cell(3,1).Shape.TextFrame.TextRange.Paragraphs(1,1).Runs(1,1).Text = `
    "Nachhaltig gewonnen: Da wir uns der Nachhaltigkeit verpflichtet haben, beziehen wir unsere Zutaten von kleinen Bauernhöfen, die ökologische Landwirtschaft betreiben. So garantieren wir nicht nur beste Qualität, sondern tragen auch zum Wohlergehen unseres Planeten bei."

$tbl6.Cell(3,2).Shape.TextFrame.TextRange.Paragraphs(1,1).Runs(1,1).Text = `
    "Elegante Verpackung: Mystic Spice Chai Tea wird in einer wunderschönen, umweltfreundlichen Verpackung geliefert, die ihn zu einem idealen Geschenk für Teeliebhaber oder zu einem luxuriösen Genuss für Sie selbst macht."

$tbl6.Cell(4,1).Shape.TextFrame.TextRange.Paragraphs(1,1).Runs(1,1).Text = `
    "Kundenzufriedenheitsgarantie: Wir stehen hinter unserem Produkt und bieten eine Zufriedenheitsgarantie."

$tbl6.Cell(4,2).Shape.TextFrame.TextRange.Paragraphs(1,1).Runs(1,1).Text = `
    "Ideal für: Teeliebhaber, gesundheitsbewusste Menschen, Liebhaber von warmen, würzigen Getränken und alle, die den reichen Geschmack des traditionellen indischen Chai entdecken möchten."
